$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the mislabeled "Grade 8" rows -> these students are actually graduates.
#    (Once no cell references the shared string "Grade 8" any more, it is dropped
#    from the workbook automatically on save.)
$ws.Range("C124:C131").Value = "Graduate"

# 2) Add the two new tracking columns: ONBOARDING YEAR (E) / ONBOARDING TERM (F).
#    Write the headers first, then copy the formatting from the existing header
#    cell (D1) so the new header cells pick up the same style as the rest of row 1.
$ws.Range("E1").Value = "ONBOARDING YEAR"
$ws.Range("F1").Value = "ONBOARDING TERM"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Populate the onboarding year/term for every student row.
$ws.Range("E2:E131").Value = 2022
$ws.Range("F2:F131").Value = 3

# 4) Size the new columns the way they were sized in the edited workbook.
$ws.Columns.Item(5).ColumnWidth = 18.736979166666668
$ws.Columns.Item(6).ColumnWidth = 23.166666666666668

# 5) Update the view: scroll down a row further and select the new onboarding cell.
$ws.Range("A113").Select() | Out-Null
$ws.Range("G125").Select() | Out-Null
